# The backwardElimination workbook holds one statsmodels OLS summary dump
# per sheet (in cell B2). The model was re-run the next day, so every
# summary's "Date:" / "Time:" header line needs to reflect the new run.
$wb = $excel.ActiveWorkbook

$oldDate = "Sat, 28 Dec 2019"
$newDate = "Sun, 29 Dec 2019"
$oldTime = "21:00:04"
$newTime = "16:11:38"

for ($i = 1; $i -le $wb.Worksheets.Count; $i++) {
    $ws = $wb.Worksheets.Item($i)
    $cell = $ws.Range("B2")
    $text = $cell.Text
    if ($text -and ($text.Contains($oldDate) -or $text.Contains($oldTime))) {
        $updated = $text.Replace($oldDate, $newDate).Replace($oldTime, $newTime)
        $cell.Value = $updated
    }
}
